$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultC")
$ws.Range("A1").Value = "What is the color of the car?  I mean, if you are looking at the car from really far away, it looks black, but I don't think it is really black.  Can you see closer up?   What is its real color?"
$ws.Activate()
$ws.Range("C24").Select()
